$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename Sex column values: "Male" -> "M", "Female" -> "F"
for ($r = 2; $r -le 35; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $v = $cell.Value2
    if ($v -eq "Male") {
        $cell.Value2 = "M"
    } elseif ($v -eq "Female") {
        $cell.Value2 = "F"
    }
}

# Rename header labels to shorter forms
$ws.Range("F1").Value2 = "Conc"
$ws.Range("E1").Value2 = "Dilut"

# Update column G width (target OOXML width 20.7109375 characters; the host
# quantizes ColumnWidth to a 1/12-character grid before padding, so feed the
# pre-image that lands closest to the desired stored width)
$ws.Columns.Item(7).ColumnWidth = 19.833333333333332

# Update the active selection
$ws.Range("O15").Select()
